# update code tinh luong
# Sheet "Lương": restructure the salary breakdown rows -- CẦN THƠ gets its own
# "Tổng công" / "Lương công tác" rows (replacing the old "Ngày công" / "Phụ cấp"
# rows), LONG XUYÊN gains a "Tổng công" / "Lương công tác" pair, SÓC TRĂNG gains
# a "Tổng công" row, and the "Tổng lương tại ..." figures at the bottom are
# refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# --- Rows 2-3: relabel "Ngày công" / "Phụ cấp" as the CẦN THƠ totals ---
$ws.Range("A2").Value = "Tổng công tại CẦN THƠ"
$ws.Range("B2").Value = 0
$ws.Range("A3").Value = "Lương công tác tại CẦN THƠ"
$ws.Range("B3").Value = 0

# --- Insert the LONG XUYÊN "Tổng công" / "Lương công tác" pair before the
#     existing "Lương cơ bản tại LONG XUYÊN" row (currently row 12) ---
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "Tổng công tại LONG XUYÊN"
$ws.Range("B12").Value = 0
$ws.Range("A13").Value = "Lương công tác tại LONG XUYÊN"
$ws.Range("B13").Value = 0

# --- Insert the SÓC TRĂNG "Tổng công" row before the existing
#     "Lương cơ bản tại SÓC TRĂNG" row (now shifted to row 22) ---
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "Tổng công tại SÓC TRĂNG"
$ws.Range("B22").Value = 26

# --- Refresh the trailing "Tổng lương tại ..." / "Tổng lương" totals, now at
#     rows 31-34 after the three inserted rows above ---
$ws.Range("B33").Value = -3500000
$ws.Range("B34").Value = -3500000
